# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) for affected leve rows
# across each crafting-class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 145898.31
$ws.Range("I132").Value = 3157.1833
$ws.Range("J132").Value = 1002345.1
$ws.Range("K132").Value = 9471.5499
$ws.Range("L132").Value = 3007035.3
$ws.Range("M132").Value = -6941.5499
$ws.Range("N132").Value = -3012095.3
$ws.Range("H141").Value = 2120
$ws.Range("I141").Value = 2356.875
$ws.Range("J141").Value = 1488.3334
$ws.Range("K141").Value = 7070.625
$ws.Range("L141").Value = 4465.0002
$ws.Range("M141").Value = -1890.625
$ws.Range("N141").Value = -14825.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1057.0483
$ws.Range("I74").Value = 888.1754
$ws.Range("K74").Value = 888.1754
$ws.Range("M74").Value = -14.17539999999997
$ws.Range("H77").Value = 1057.0483
$ws.Range("I77").Value = 888.1754
$ws.Range("K77").Value = 4440.876999999999
$ws.Range("M77").Value = -72.8769999999995
$ws.Range("H102").Value = 2400
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = $null
$ws.Range("H132").Value = 181105.48
$ws.Range("I132").Value = 7267.892
$ws.Range("J132").Value = 502705.06
$ws.Range("K132").Value = 21803.676
$ws.Range("L132").Value = 1508115.18
$ws.Range("M132").Value = -19273.676
$ws.Range("N132").Value = -1513175.18

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2300
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24721.49
$ws.Range("I31").Value = 27248.23
$ws.Range("J31").Value = 18562.562
$ws.Range("K31").Value = 27248.23
$ws.Range("L31").Value = 18562.562
$ws.Range("M31").Value = -26953.23
$ws.Range("N31").Value = -19152.562
$ws.Range("H34").Value = 24721.49
$ws.Range("I34").Value = 27248.23
$ws.Range("J34").Value = 18562.562
$ws.Range("K34").Value = 27248.23
$ws.Range("L34").Value = 18562.562
$ws.Range("M34").Value = -27046.23
$ws.Range("N34").Value = -18966.562
$ws.Range("H62").Value = 3281.9333
$ws.Range("I62").Value = 2484.5454
$ws.Range("J62").Value = 5474.75
$ws.Range("K62").Value = 2484.5454
$ws.Range("L62").Value = 5474.75
$ws.Range("M62").Value = -1860.5454
$ws.Range("N62").Value = -6722.75
$ws.Range("H65").Value = 3281.9333
$ws.Range("I65").Value = 2484.5454
$ws.Range("J65").Value = 5474.75
$ws.Range("K65").Value = 12422.727
$ws.Range("L65").Value = 27373.75
$ws.Range("M65").Value = -9302.726999999999
$ws.Range("N65").Value = -33613.75
$ws.Range("H86").Value = 48078720
$ws.Range("I86").Value = 62501716
$ws.Range("J86").Value = 2056.6667
$ws.Range("K86").Value = 62501716
$ws.Range("L86").Value = 2056.6667
$ws.Range("M86").Value = -62500593
$ws.Range("N86").Value = -4302.6667
$ws.Range("H89").Value = 48078720
$ws.Range("I89").Value = 62501716
$ws.Range("J89").Value = 2056.6667
$ws.Range("K89").Value = 312508580
$ws.Range("L89").Value = 10283.3335
$ws.Range("M89").Value = -312502964
$ws.Range("N89").Value = -21515.3335
$ws.Range("H122").Value = 787
$ws.Range("I122").Value = 878.6667
$ws.Range("J122").Value = 512
$ws.Range("K122").Value = 2636.0001
$ws.Range("L122").Value = 1536
$ws.Range("M122").Value = -186.0001000000002
$ws.Range("N122").Value = -6436
$ws.Range("H132").Value = 24062.705
$ws.Range("I132").Value = 34369.535
$ws.Range("J132").Value = 1976.6428
$ws.Range("K132").Value = 103108.605
$ws.Range("L132").Value = 5929.928400000001
$ws.Range("M132").Value = -100578.605
$ws.Range("N132").Value = -10989.9284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1768.1
$ws.Range("I5").Value = 1623.7142
$ws.Range("J5").Value = 1798.7273
$ws.Range("K5").Value = 4871.142599999999
$ws.Range("L5").Value = 5396.1819
$ws.Range("M5").Value = -4759.142599999999
$ws.Range("N5").Value = -5620.1819
$ws.Range("H68").Value = 333634
$ws.Range("I68").Value = 500151
$ws.Range("J68").Value = 600
$ws.Range("K68").Value = 1500453
$ws.Range("L68").Value = 1800
$ws.Range("M68").Value = -1499642
$ws.Range("N68").Value = -3422
$ws.Range("H71").Value = 333634
$ws.Range("I71").Value = 500151
$ws.Range("J71").Value = 600
$ws.Range("K71").Value = 4501359
$ws.Range("L71").Value = 5400
$ws.Range("M71").Value = -4497303
$ws.Range("N71").Value = -13512
$ws.Range("H131").Value = 186059.27
$ws.Range("I131").Value = 447.27274
$ws.Range("J131").Value = 233541.39
$ws.Range("K131").Value = 1341.81822
$ws.Range("L131").Value = 700624.17
$ws.Range("M131").Value = 3698.18178
$ws.Range("N131").Value = -710704.17
$ws.Range("H135").Value = 1768.1
$ws.Range("I135").Value = 1623.7142
$ws.Range("J135").Value = 1798.7273
$ws.Range("K135").Value = 14613.4278
$ws.Range("L135").Value = 16188.5457
$ws.Range("M135").Value = -12078.4278
$ws.Range("N135").Value = -21258.5457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1819.5
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 1926
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 5778
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -10718
$ws.Range("H132").Value = 32813.875
$ws.Range("I132").Value = 976.6
$ws.Range("J132").Value = 85876
$ws.Range("K132").Value = 2929.8
$ws.Range("L132").Value = 257628
$ws.Range("M132").Value = -399.8000000000002
$ws.Range("N132").Value = -262688
$ws.Range("H134").Value = 31800
$ws.Range("J134").Value = 31800
$ws.Range("L134").Value = 95400
$ws.Range("N134").Value = -100470
$ws.Range("H136").Value = 26000
$ws.Range("J136").Value = 26000
$ws.Range("L136").Value = 78000
$ws.Range("N136").Value = -83100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 48000
$ws.Range("J98").Value = 48000
$ws.Range("L98").Value = 48000
$ws.Range("N98").Value = -53990
$ws.Range("H132").Value = 259257.55
$ws.Range("I132").Value = 77762.15
$ws.Range("J132").Value = 504276.34
$ws.Range("K132").Value = 233286.45
$ws.Range("L132").Value = 1512829.02
$ws.Range("M132").Value = -230756.45
$ws.Range("N132").Value = -1517889.02
$ws.Range("H140").Value = 45822.9
$ws.Range("J140").Value = 45822.9
$ws.Range("L140").Value = 45822.9
$ws.Range("N140").Value = -56182.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3440.0444
$ws.Range("I132").Value = 705.4828
$ws.Range("J132").Value = 8396.4375
$ws.Range("K132").Value = 2116.4484
$ws.Range("L132").Value = 25189.3125
$ws.Range("M132").Value = 413.5515999999998
$ws.Range("N132").Value = -30249.3125
$ws.Range("H141").Value = 50747.918
$ws.Range("J141").Value = 50747.918
$ws.Range("L141").Value = 50747.918
$ws.Range("N141").Value = -61107.918

